$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 883.8461
$ws.Range("I9").Value = 298.75
$ws.Range("K9").Value = 298.75
$ws.Range("M9").Value = -129.75
# Row 18
$ws.Range("H18").Value = 250
$ws.Range("I18").Value = 250
$ws.Range("K18").Value = 250
$ws.Range("M18").Value = 34
# Row 32
$ws.Range("H32").Value = 1632.3334
$ws.Range("I32").Value = 1749.75
$ws.Range("K32").Value = 1749.75
$ws.Range("M32").Value = -1423.75
# Row 43
$ws.Range("H43").Value = 4449.5713
$ws.Range("I43").Value = 4200.0557
$ws.Range("K43").Value = 4200.0557
$ws.Range("M43").Value = -4131.0557
# Row 86
$ws.Range("H86").Value = 7403
$ws.Range("I86").Value = 6859.8
$ws.Range("K86").Value = 6859.8
$ws.Range("M86").Value = -5736.8
# Row 89
$ws.Range("H89").Value = 7403
$ws.Range("I89").Value = 6859.8
$ws.Range("K89").Value = 34299
$ws.Range("M89").Value = -28683
# Row 113
$ws.Range("H113").Value = 8820.1
$ws.Range("I113").Value = 8524.5
$ws.Range("K113").Value = 8524.5
$ws.Range("M113").Value = -5270.5
# Row 116
$ws.Range("H116").Value = 2167.5
$ws.Range("I116").Value = 2101
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 2101
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = 1341
$ws.Range("N116").Value = -9384
# Row 132
$ws.Range("H132").Value = 2297.8462
$ws.Range("I132").Value = 2021.619
$ws.Range("K132").Value = 6064.857
$ws.Range("M132").Value = -3534.857
# Row 135
$ws.Range("H135").Value = 500
$ws.Range("I135").Value = 500
$ws.Range("K135").Value = 4500
$ws.Range("M135").Value = -1965
# Row 138
$ws.Range("H138").Value = 4753.1333
$ws.Range("J138").Value = 4500
$ws.Range("L138").Value = 13500
$ws.Range("N138").Value = -23780

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4094.7368
$ws.Range("I45").Value = 3000
$ws.Range("K45").Value = 3000
$ws.Range("M45").Value = -2623
# Row 63
$ws.Range("H63").Value = 2471
$ws.Range("I63").Value = 2310.75
$ws.Range("J63").Value = 2898.3333
$ws.Range("K63").Value = 2310.75
$ws.Range("L63").Value = 2898.3333
$ws.Range("M63").Value = -1624.75
$ws.Range("N63").Value = -4270.3333
# Row 66
$ws.Range("H66").Value = 2471
$ws.Range("I66").Value = 2310.75
$ws.Range("J66").Value = 2898.3333
$ws.Range("K66").Value = 11553.75
$ws.Range("L66").Value = 14491.6665
$ws.Range("M66").Value = -8121.75
$ws.Range("N66").Value = -21355.6665
# Row 102
$ws.Range("H102").Value = 3063.2727
$ws.Range("J102").Value = 1500
$ws.Range("L102").Value = 1500
$ws.Range("N102").Value = -4744
# Row 122
$ws.Range("H122").Value = 1914.7858
$ws.Range("I122").Value = 1014.44446
$ws.Range("J122").Value = 3535.4
$ws.Range("K122").Value = 3043.33338
$ws.Range("L122").Value = 10606.2
$ws.Range("M122").Value = -593.33338
$ws.Range("N122").Value = -15506.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 3430
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = ""
# Row 80
$ws.Range("H80").Value = 478.42856
$ws.Range("I80").Value = 520
$ws.Range("J80").Value = 374.5
$ws.Range("K80").Value = 520
$ws.Range("L80").Value = 374.5
$ws.Range("M80").Value = 478
$ws.Range("N80").Value = -2370.5
# Row 83
$ws.Range("H83").Value = 478.42856
$ws.Range("I83").Value = 520
$ws.Range("J83").Value = 374.5
$ws.Range("K83").Value = 2600
$ws.Range("L83").Value = 1872.5
$ws.Range("M83").Value = 2392
$ws.Range("N83").Value = -11856.5
# Row 86
$ws.Range("H86").Value = 3143.2856
$ws.Range("I86").Value = 2689.7778
$ws.Range("J86").Value = 3959.6
$ws.Range("K86").Value = 2689.7778
$ws.Range("L86").Value = 3959.6
$ws.Range("M86").Value = -1566.7778
$ws.Range("N86").Value = -6205.6
# Row 89
$ws.Range("H89").Value = 3143.2856
$ws.Range("I89").Value = 2689.7778
$ws.Range("J89").Value = 3959.6
$ws.Range("K89").Value = 13448.889
$ws.Range("L89").Value = 19798
$ws.Range("M89").Value = -7832.888999999999
$ws.Range("N89").Value = -31030
# Row 99
$ws.Range("H99").Value = 2641.2144
$ws.Range("I99").Value = 2179.818
$ws.Range("K99").Value = 2179.818
$ws.Range("M99").Value = -681.8180000000002
# Row 134
$ws.Range("H134").Value = 2248.3333
$ws.Range("I134").Value = 2248.3333
$ws.Range("K134").Value = 6744.999899999999
$ws.Range("M134").Value = -4209.999899999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 2995
$ws.Range("I41").Value = 2995
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2995
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2567
$ws.Range("N41").Value = ""
# Row 107
$ws.Range("H107").Value = 1175
$ws.Range("I107").Value = 450
$ws.Range("K107").Value = 450
$ws.Range("M107").Value = 1470
# Row 122
$ws.Range("H122").Value = 2417.5
$ws.Range("I122").Value = 2242.2222
$ws.Range("K122").Value = 6726.6666
$ws.Range("M122").Value = -4276.6666
# Row 132
$ws.Range("H132").Value = 2389
$ws.Range("I132").Value = 2390.1155
$ws.Range("J132").Value = 2374.5
$ws.Range("K132").Value = 7170.3465
$ws.Range("L132").Value = 7123.5
$ws.Range("M132").Value = -4640.3465
$ws.Range("N132").Value = -12183.5
# Row 134
$ws.Range("H134").Value = 1706.4375
$ws.Range("I134").Value = 1486.8667
$ws.Range("K134").Value = 4460.6001
$ws.Range("M134").Value = -1925.6001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 45
$ws.Range("H45").Value = 1933
$ws.Range("J45").Value = 1933
$ws.Range("L45").Value = 5799
$ws.Range("N45").Value = -6863
# Row 103
$ws.Range("H103").Value = 31765.5
$ws.Range("J103").Value = 1399.5
$ws.Range("L103").Value = 4198.5
$ws.Range("N103").Value = -5956.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 23
$ws.Range("H23").Value = 3420.3333
$ws.Range("I23").Value = 2001
$ws.Range("J23").Value = 3704.2
$ws.Range("K23").Value = 2001
$ws.Range("L23").Value = 3704.2
$ws.Range("M23").Value = -1778
$ws.Range("N23").Value = -4150.2
# Row 80
$ws.Range("H80").Value = 5384.2856
$ws.Range("I80").Value = 2668.25
$ws.Range("K80").Value = 2668.25
$ws.Range("M80").Value = -1670.25
# Row 83
$ws.Range("H83").Value = 5384.2856
$ws.Range("I83").Value = 2668.25
$ws.Range("K83").Value = 13341.25
$ws.Range("M83").Value = -8349.25
# Row 102
$ws.Range("H102").Value = 1255.6
$ws.Range("I102").Value = 500
$ws.Range("J102").Value = 1444.5
$ws.Range("K102").Value = 500
$ws.Range("L102").Value = 1444.5
$ws.Range("M102").Value = 1122
$ws.Range("N102").Value = -4688.5
# Row 122
$ws.Range("H122").Value = 3545.4614
$ws.Range("I122").Value = 3553
$ws.Range("K122").Value = 10659
$ws.Range("M122").Value = -8209
# Row 126
$ws.Range("H126").Value = 3996.4285
$ws.Range("I126").Value = 3996.25
$ws.Range("J126").Value = 3996.6667
$ws.Range("K126").Value = 11988.75
$ws.Range("L126").Value = 11990.0001
$ws.Range("M126").Value = -9518.75
$ws.Range("N126").Value = -16930.0001
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = ""

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 2990
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2990
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2990
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = -3270
# Row 46
$ws.Range("H46").Value = 4483.8213
$ws.Range("I46").Value = 4572.727
$ws.Range("J46").Value = 4426.294
$ws.Range("K46").Value = 4572.727
$ws.Range("L46").Value = 4426.294
$ws.Range("M46").Value = -4384.727
$ws.Range("N46").Value = -4802.294
# Row 54
$ws.Range("H54").Value = 45042
$ws.Range("I54").Value = 35000
$ws.Range("J54").Value = 55084
$ws.Range("K54").Value = 35000
$ws.Range("L54").Value = 55084
$ws.Range("M54").Value = -34356
$ws.Range("N54").Value = -56372
# Row 82
$ws.Range("H82").Value = 950.1667
$ws.Range("I82").Value = 1100.6666
$ws.Range("J82").Value = 799.6667
$ws.Range("K82").Value = 1100.6666
$ws.Range("L82").Value = 799.6667
$ws.Range("M82").Value = -739.6666
$ws.Range("N82").Value = -1521.6667
# Row 85
$ws.Range("H85").Value = 950.1667
$ws.Range("I85").Value = 1100.6666
$ws.Range("J85").Value = 799.6667
$ws.Range("K85").Value = 1100.6666
$ws.Range("L85").Value = 799.6667
$ws.Range("M85").Value = 147.3334
$ws.Range("N85").Value = -3295.6667
# Row 136
$ws.Range("H136").Value = 4610
$ws.Range("I136").Value = 4610
$ws.Range("K136").Value = 13830
$ws.Range("M136").Value = -11280

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3125.75
$ws.Range("I62").Value = 2833.3333
$ws.Range("K62").Value = 2833.3333
$ws.Range("M62").Value = -2209.3333
# Row 65
$ws.Range("H65").Value = 3125.75
$ws.Range("I65").Value = 2833.3333
$ws.Range("K65").Value = 14166.6665
$ws.Range("M65").Value = -11046.6665
# Row 81
$ws.Range("H81").Value = 6999
$ws.Range("J81").Value = 5000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -12122
# Row 84
$ws.Range("H84").Value = 6999
$ws.Range("J84").Value = 5000
$ws.Range("L84").Value = 50000
$ws.Range("N84").Value = -60608
# Row 107
$ws.Range("H107").Value = 433.2857
$ws.Range("I107").Value = 397.25
$ws.Range("K107").Value = 1191.75
$ws.Range("M107").Value = 728.25
# Row 136
$ws.Range("H136").Value = 755
$ws.Range("J136").Value = 1500
$ws.Range("L136").Value = 4500
$ws.Range("N136").Value = -9600
# Row 140
$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360
# Row 141
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360
